$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds plain text in the source data (values like
# "1.00", "10.20", "0.0740" must keep their exact digits/trailing
# zeros). Pre-format the cells that are about to receive new numeric-
# looking text as Text ("@") so Excel does not silently convert them
# to real numbers (which would drop formatting like trailing zeros).
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D12").NumberFormat = "@"
$ws.Range("D14:D25").NumberFormat = "@"
$ws.Range("D27:D28").NumberFormat = "@"
$ws.Range("D30:D32").NumberFormat = "@"
$ws.Range("D34:D36").NumberFormat = "@"
$ws.Range("D39:D40").NumberFormat = "@"
$ws.Range("D42:D48").NumberFormat = "@"
$ws.Range("D50:D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.447.01"
$ws.Range("E2").Value = "  +3.22%  "
$ws.Range("D3").Value = "2.252.51"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "302.73"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").Value = "90.84"
$ws.Range("E6").Value = "  +3.48%  "
$ws.Range("D7").Value = "0.527"
$ws.Range("E7").Value = "  +2.71%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.476"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").Value = "31.91"
$ws.Range("E10").Value = "  +3.33%  "
$ws.Range("D11").Value = "52.69"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "0.0796"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "6.54"
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("D15").Value = "2.607.08"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("D16").Value = "14.12"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "2.248.14"
$ws.Range("E17").Value = "  +3.55%  "
$ws.Range("D18").Value = "0.756"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("D19").Value = "41.425.06"
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("D20").Value = "12.21"
$ws.Range("E20").Value = "  +8.18%  "
$ws.Range("D21").Value = "0.0₃0899"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("D22").Value = "5.86"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("D23").Value = "66.45"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("D24").Value = "239.61"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "1.89"
$ws.Range("E27").Value = "  +4.71%  "
$ws.Range("D28").Value = "23.83"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("D30").Value = "9.44"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").Value = "159.83"
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("D32").Value = "33.96"
$ws.Range("E32").Value = "  +5.89%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "5.13"
$ws.Range("E34").Value = "  +3.30%  "
$ws.Range("D35").Value = "0.0738"
$ws.Range("E35").Value = "  +3.39%  "
$ws.Range("D36").Value = "2.99"
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("D39").Value = "16.51"
$ws.Range("E39").Value = "  +5.37%  "
$ws.Range("D40").Value = "0.103"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("D42").Value = "3.88"
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("D43").Value = "2.048.70"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "19.49"
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").Value = "0.0276"
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("D46").Value = "10.22"
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.04"
$ws.Range("E47").Value = "  +6.57%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "2.83"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("E49").Value = "  +3.51%  "
$ws.Range("D50").Value = "72.17"
$ws.Range("E50").Value = "  +6.23%  "
$ws.Range("D51").Value = "1.14"
$ws.Range("E51").Value = "  +1.43%  "
